# Updated symbol list on Sat Jan 28 14:37:01 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) quotes for the
# crypto-ranking rows on the active sheet. Values are stored as plain text
# (matching the source feed's inlineStr cells), so each cell is forced to
# Text format before the write and reset to the default "Normal" style
# afterwards so no stray number-format style lingers on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2";  Value = "305.92" }
    @{ Cell = "E2";  Value = "0.82%" }

    @{ Cell = "D3";  Value = "38.24" }
    @{ Cell = "E3";  Value = "7.00%" }

    @{ Cell = "E4";  Value = "1.11%" }

    @{ Cell = "D5";  Value = "0.08064" }
    @{ Cell = "E5";  Value = "1.09%" }

    @{ Cell = "D6";  Value = "1.940" }
    @{ Cell = "E6";  Value = "4.51%" }

    @{ Cell = "E7";  Value = "1.49%" }

    @{ Cell = "D8";  Value = "7.949" }
    @{ Cell = "E8";  Value = "2.42%" }

    @{ Cell = "D9";  Value = "0.9292" }
    @{ Cell = "E9";  Value = "0.71%" }

    @{ Cell = "D10"; Value = "0.1442" }
    @{ Cell = "E10"; Value = "13.35%" }

    @{ Cell = "D11"; Value = "0.1920" }
    @{ Cell = "E11"; Value = "2.08%" }

    @{ Cell = "D12"; Value = "0.09025" }
    @{ Cell = "E12"; Value = "0.95%" }

    @{ Cell = "D13"; Value = "0.03514" }
    @{ Cell = "E13"; Value = "2.69%" }

    @{ Cell = "D14"; Value = "0.09788" }
    @{ Cell = "E14"; Value = "-0.57%" }

    @{ Cell = "D15"; Value = "0.001391" }
    @{ Cell = "E15"; Value = "-1.18%" }

    @{ Cell = "D16"; Value = "0.006163" }
    @{ Cell = "E16"; Value = "-2.31%" }

    @{ Cell = "E17"; Value = "-3.62%" }

    @{ Cell = "E18"; Value = "3.65%" }

    @{ Cell = "E19"; Value = "1.65%" }

    @{ Cell = "D20"; Value = "0.1311" }
    @{ Cell = "E20"; Value = "-2.18%" }

    @{ Cell = "D21"; Value = "4.788" }
    @{ Cell = "E21"; Value = "-0.32%" }

    @{ Cell = "D22"; Value = "0.2405" }
    @{ Cell = "E22"; Value = "2.67%" }

    @{ Cell = "D23"; Value = "0.04350" }
    @{ Cell = "E23"; Value = "-0.12%" }

    @{ Cell = "D24"; Value = "0.001232" }
    @{ Cell = "E24"; Value = "-0.32%" }

    @{ Cell = "D25"; Value = "0.004119" }
    @{ Cell = "E25"; Value = "-14.94%" }

    @{ Cell = "E27"; Value = "-0.04%" }

    @{ Cell = "D39"; Value = "0.02073" }
    @{ Cell = "E39"; Value = "8.06%" }

    @{ Cell = "D40"; Value = "0.05029" }
    @{ Cell = "E40"; Value = "-1.43%" }

    @{ Cell = "D41"; Value = "0.007474" }
    @{ Cell = "E41"; Value = "-1.20%" }

    @{ Cell = "E42"; Value = "-0.39%" }

    @{ Cell = "D43"; Value = "0.1347" }
    @{ Cell = "E43"; Value = "0.34%" }

    @{ Cell = "E44"; Value = "1.38%" }

    @{ Cell = "D45"; Value = "0.008918" }
    @{ Cell = "E45"; Value = "-9.64%" }

    @{ Cell = "D46"; Value = "0.00006183" }
    @{ Cell = "E46"; Value = "-0.47%" }

    @{ Cell = "E47"; Value = "-0.05%" }

    @{ Cell = "D48"; Value = "0.002810" }

    @{ Cell = "E49"; Value = "27.85%" }

    @{ Cell = "E50"; Value = "-0.05%" }

    @{ Cell = "E51"; Value = "-0.05%" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    # Force Text format so the numeric-looking string isn't auto-converted
    # to a number, then write the value, then drop back to the default
    # "Normal" style so no stray number-format style is left on the cell.
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
